# cv124101a.xlsx - "correção nos dados e inicio da analise PNAD 2009"
#
# 1) Row 2, column B held the stray pandas artifact label
#    "unnamed: 1_level_1" - replace it with the real column header "total".
# 2) Several rows in the sheet were pure category captions with no data
#    underneath them (sexo / cor ou raça / grupos de idade / nível de
#    instrução / classes de rendimento mensal domiciliar per capita) plus
#    two trailing footnote-only rows at the very end. These are removed,
#    which shifts every data row below them upward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "total"

# Delete from the bottom up so earlier row numbers stay valid as we go.
$rowsToDelete = @(35, 34, 27, 19, 13, 8, 5)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}
